# Append a 4th row of formulas, mirroring the pattern already used in rows 1-3:
# col A holds a standalone formula, and cols B:C hold a shared formula group.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Formula = "=E4*F4"
$ws.Range("B4:C4").Formula = "=F4*G4"

# Move the selection back to A1 (the diff drops the explicit <selection> that
# pointed at A3, i.e. the view resets to the sheet's default top-left cell).
$ws.Range("A1").Select()
